# day4 notes update:
#  1. Footer "datetimeFigureOut" field text 9/17/2016 -> 9/18/2016
#     (lives on the Slide Master + every Custom Layout's Date placeholder)
#  2. Slide 12: typo fix "in conjunctions" -> "in conjunction"
#  3. Slide 14: rewrite the VBA paragraph to a shorter sentence and add
#     two new paragraphs ("(... I had more here.)" plus a blank line)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the date placeholder everywhere it is cached: slide master
#    and all custom (slide) layouts.
# ---------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "9/17/2016") {
                $tr.Text = "9/18/2016"
            }
        }
    }
}

Update-DateShape $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 12 - fix "in conjunctions" -> "in conjunction"
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(2)
$tr12 = $sh12.TextFrame.TextRange
$full12 = $tr12.Text
$marker = "and use it to edit code"
$idx12 = $full12.IndexOf($marker)
$len12 = $full12.Length - $idx12
$sub12 = $tr12.Characters($idx12 + 1, $len12)
$sub12.Text = "and use it to edit code in many languages, in conjunction with separate compilers, interpreters, etc."

# ---------------------------------------------------------------------
# 3) Slide 14 - shrink the VBA rant and add the new follow-up paragraphs
# ---------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(2)
$tr14 = $sh14.TextFrame.TextRange

$ellipsis = [char]0x2026

# Replace the long VBA sentence (3 runs) with a single shorter run.
$full14 = $tr14.Text
$startMarker = "VBA is a condescending"
$endMarker = "each design decision."
$sIdx = $full14.IndexOf($startMarker)
$eIdx = $full14.IndexOf($endMarker) + $endMarker.Length
$oldLen = $eIdx - $sIdx
$vbaRange = $tr14.Characters($sIdx + 1, $oldLen)
$newSentence = "VBA is a condescending, under-designed, overcomplicated dumpster fire of a language" + $ellipsis
$vbaRange.Text = $newSentence

# Insert the two new paragraphs right after the (still single) VBA
# paragraph, before the pre-existing blank paragraph / closing remark.
$full14b = $tr14.Text
$insertAt = $sIdx + $newSentence.Length + 1
$insertion = "`r(" + $ellipsis + " I had more here.)`r"
$tailRange = $tr14.Characters($insertAt, $full14b.Length - $insertAt + 1)
$tailRange.Text = $insertion + $tailRange.Text

# Normalize the shape's autofit so it no longer carries the stale
# font-scale cache from the much longer paragraph.
$sh14.TextFrame.AutoSize = 0
$sh14.TextFrame.AutoSize = 2
